$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.078810386556763298
$ws.Range("B1").Value = 0.078810385301768662
$ws.Range("A2").Value = 0.019229882981893553
$ws.Range("B2").Value = -0.019229884278862252
$ws.Range("A3").Value = 0.031969756374199597
$ws.Range("B3").Value = -0.031969757616886509
$ws.Range("A4").Value = 0.00055431271148150072
$ws.Range("B4").Value = -0.0005543140649239978
